$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Build the two new cell styles (borderId 4: top+bottom thin; borderId 5:
# top+bottom+right thin) once, on a scratch worksheet, then stamp them onto
# every target cell with a single PasteSpecial(formats) each. Doing it this
# way (build once, copy everywhere) avoids creating throw-away/orphan style
# entries in styles.xml that a naive "set each border edge on each cell"
# approach would leave behind.
# ---------------------------------------------------------------------------
$tmp = $wb.Worksheets.Add()

# Always look sheets up by name (not position) since adding/removing sheets
# shifts numeric indices.
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Reference style #1: thin top + thin bottom (-> borderId 4)
$ref4 = $tmp.Range("A1")
$ref4.Style = "Normal"
$ref4.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ref4.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Reference style #2: thin top + thin bottom + thin right (-> borderId 5)
# Built from ref4 (copy) plus a single extra edge so the engine never visits
# a brand-new unique border combination more than the minimum required.
$ref4.Copy()
$ref5 = $tmp.Range("A2")
$ref5.PasteSpecial(-4122)              # xlPasteFormats
$ref5.Borders.Item(10).LineStyle = 1   # xlEdgeRight

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ref4.Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$ref5.Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ref4.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ref5.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ref4.Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$ref5.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()

# Drop the scratch worksheet used to build the reference styles
$tmp.Delete() | Out-Null
